$d = $word.ActiveDocument

# The site footer boilerplate ("Ver no Jupiter Salvar em pdf Salvar em docx"
# and the "© 2020 ... Jekyll and Github pages ..." copyright line), together
# with the blank paragraph that separates them from the preceding
# "Requisitos" content, was dropped from the generated page on this site
# rebuild. Locate the boilerplate by its distinctive text rather than a
# hard-coded paragraph index, then remove it along with its neighbouring
# blank paragraph, leaving the trailing blank paragraph (the one right
# before the final page-break paragraph) untouched.

$marker = "Ver no Jupiter"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*$marker*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # Delete highest index first so lower indices remain valid:
    #   targetIndex+1 -> copyright paragraph
    #   targetIndex   -> "Ver no Jupiter Salvar em pdf Salvar em docx"
    #   targetIndex-1 -> blank paragraph right before it
    $d.Paragraphs.Item($targetIndex + 1).Range.Delete()
    $d.Paragraphs.Item($targetIndex).Range.Delete()
    $d.Paragraphs.Item($targetIndex - 1).Range.Delete()
}
